$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Craciun")

# Mark rows 10-14 (Reserved / bought column E) as "Y"
$ws.Range("E10").Value = "Y"
$ws.Range("E11").Value = "Y"
$ws.Range("E12").Value = "Y"
$ws.Range("E13").Value = "Y"
$ws.Range("E14").Value = "Y"

# Update puzzle age range text
$ws.Range("A17").Value = "Puzzle (orice 3-5 ani)"

# New wishlist rows 19-25 -- written in the original authoring order so that
# the shared-string table comes out in the same order as the source file.

# Row 19 - LEGO Princess Market 43246
$ws.Range("A19").Value = "LEGO Princess Market 43246"
$ws.Range("C19").Value = "https://www.amazon.de/-/en/LEGO-Disney-Princess-Adventure-Cinderella/dp/B0CFW17427/ref=sr_1_1?crid=2S57I7DOW3RGJ&dib=eyJ2IjoiMSJ9.O4IX-zfr64af8opy3QSNF2ZT80sQLR4S5LkCENoLOZUsOdoigdNbq6JzcaiqiIhoBWzcD-sLcS7kUiCmkJsR2DBrUBl8Vj8KjQf9SRmP70sim25K3BZKQhNOWUmkrZnx49L9v6SHlL4TpdyPLWw9r9rGgq_Cv-AVhS4CqKj2YZyb2Yj5ONFOsQabxBcvsWRtRlgX_69FZPJg8Cn2zAdNA5jMwnxrT34Np7QsJi_NcHR_9ni-cyHdC-08q-if4qriU2xAiByK2HzpB3eBMbD9iQqD9J8RosV8sq6L4UKlS5A.y54hzRPnToDATsiTHnfK_onj2MxmliLgpehHCQBg9WE&dib_tag=se&keywords=lego+43246&qid=1729516242&sprefix=lego+43246%2Caps%2C201&sr=8-1"
$ws.Range("D19").Value = "68 EUR"
$ws.Range("B19").Value = "https://m.media-amazon.com/images/I/81kO40+M7UL._AC_SX679_.jpg"

# Row 20 - Catel cu telecomanda
$ws.Range("A20").Value = "Catel cu telecomanda"
$ws.Range("C20").Value = "https://www.amazon.de/-/en/furReal-Walkin-Interactive-Bouncing-Effects/dp/B0CRBJ8JFR/ref=sr_1_2?crid=1VQJRA7B0MIQV&dib=eyJ2IjoiMSJ9.53Uu6FxJ4b6G8ZdU8vF8f5_69fBrUXbiB0mmwb1j4o3FN5oTVqD_WxOQ-IU9ZCdLyRupmmDs2Bni1O65n56dwqXUKAQ_hfty5U3xd6uubP6GEKV_XpgF4J5s--YdZWRl4b9kz0DjalmySz_sZaxBs7TKoEQrpKoOIm0HBvMOyJci82LMGrfxTtnRQVckEemvsf_-Ukh732hHa8ZrScNKi2kpUIJZSlMODf42dJNSYZY.WOEefNOMKp6eIVjyUoZ4k5cSLQjkeM2JF-DkR-rmPyE&dib_tag=se&keywords=furreal+friends+hund+gogo&qid=1729516311&sprefix=furreal+gogo+hun%2Caps%2C104&sr=8-2"
$ws.Range("D20").Value = "37 EUR"
$ws.Range("B20").Value = "https://m.media-amazon.com/images/I/71A+tWCmVlL._AC_SX679_.jpg"

# Row 23 - Set potiuni
$ws.Range("A23").Value = "Set potiuni"
$ws.Range("C23").Value = "https://www.amazon.de/-/en/Mixies-Cauldron-Purple-Interactive-Conjure/dp/B0BX6YW5N9/ref=sr_1_1?crid=358IKESIBXVKO&dib=eyJ2IjoiMSJ9.oFY88I1D4MFo5A9ePZtzXjbQS4l-SpFJkwgAbWg7MMhtOcENlN6wulpw416W8xsS2fLUPGO47TnycXTPWIxHE2CDymxryzuQt0W7vBohROKY7k9Gxd5wxAniONIj00EQFSdEFqwA14QjSRc8QUbA489aooXS7DEeynFt9EmUOEY8HuoiqkTF2L2tH2U6weSLdM-zRruSR3jTRicqwIRaJ7KS33hoqDz6YFgspU9AJP0G5bDKdrcfvCTO2WeSnZDhWVx8KzBe_NmTT37oDFQSA7zNCmpXIAokJ02pSwIx80U.YJEOHI_2v0L6eI_rOrgvMkf9kcId4bmGS2Bwd1uJk7k&dib_tag=se&keywords=my+magic+mixies&qid=1729516398&sprefix=my+magic+mix%2Caps%2C99&sr=8-1"
$ws.Range("D23").Value = "66 EUR"
$ws.Range("B23").Value = "https://m.media-amazon.com/images/I/81yFvCyxfLL._AC_SX679_.jpg"

# Row 22 - Plus Stitch cu sunet si miros
$ws.Range("A22").Value = "Plus Stitch cu sunet si miros"
$ws.Range("C22").Value = "https://www.amazon.de/-/en/Stitch-Disney-Sound-Scent-Small/dp/B0C6YL2H1D/ref=sr_1_1?crid=6EMXNMGZ24N&dib=eyJ2IjoiMSJ9.GuA80LmF0yuVPEU-kqsv9uLayaq9ADpBvReu1nPc2asAtEJmgEUUTq10j-altcYb3t9k1YyGDPJevWY5WAl7q5HRvhO8xbA4PHlUbBc_VRiNiXsoIrryF9CCiw_8Ymsuuje7-xai4s4kYksI-f3wNk23ZAm99jSMVpGYeEs9Du2EbsLhExjtRQruu-kR9luH8MB5o7ECPnd7ASgVz6ZW3uqoxwIynSUbp8zjBukwaj0GZ0w34C_HnUySu_HAhnPT-3Wp0wgwUcFnsx4PPpEmbs-Hw7QzcljVB4diYLvaOs0.0Y9HuQtkgUCAVRkvEzI6CmJ_iNtgsbo1m-JlnVZ6QXk&dib_tag=se&keywords=stitch+soundfunktion+duft+15cm&qid=1729516484&sprefix=stitch+soundfunktion+duft+15cm%2Caps%2C76&sr=8-1"
$ws.Range("B22").Value = "https://m.media-amazon.com/images/I/81wmHRHoM8L._AC_SX679_.jpg"
$ws.Range("D22").Value = "17 EUR"

# Row 24 - Gravitrax Junior Starter Set L - Ice
$ws.Range("C24").Value = "https://www.amazon.de/-/en/Ravensburger-GraviTrax-Junior-Starter-Environmentally/dp/B0BSX847R3/ref=sr_1_2?crid=1X7LU8OFXMUDY&dib=eyJ2IjoiMSJ9.HacZPOrg5jUAw8zb8YX2rP9wurG1VPUH9ajjsfsttYOMungu-BuQav_M7W7LDs_HaLJiA4CCD1XsLUL7P628PHiaWR93cn8kXvmZKpxugWklhxzMUdAxO8NonF9N-uGQB4B7-4lX_3SD91JPzXuOLiPvsGDu3Menod9Pg6gkgWDHxOuh2ZRNrm5yFve3aLzysRcmg5diMQ92Yws4UuWjJvEXvD77q8MLTFxrgWAXOErjvtLv-E_gwG34fWMeduft6MSmf5F-oh8goqjjolGhi-h-dQQ5YecsoFCB3ncrJys.9EKnkXvyxCPIKZXs9Mx9U7m2t-LMTsbMy9fr4BxJK6M&dib_tag=se&keywords=gravitrax+junior+starterset&qid=1729516745&sprefix=gravitrax+junior%2Caps%2C98&sr=8-2"
$ws.Range("B24").Value = "https://m.media-amazon.com/images/I/715c7oaL5bL._AC_SX679_.jpg"
$ws.Range("A24").Value = "Gravitrax Junior Starter Set L - Ice"
$ws.Range("D24").Value = "36 EUR"

# Row 25 - Gravitrax Junior Starter Set Frozen
$ws.Range("C25").Value = "https://www.amazon.de/-/en/Ravensburger-73385-GraviTrax-Starter-Children/dp/B0CSZ63YV3/ref=sr_1_5?crid=1X7LU8OFXMUDY&dib=eyJ2IjoiMSJ9.HacZPOrg5jUAw8zb8YX2rP9wurG1VPUH9ajjsfsttYOMungu-BuQav_M7W7LDs_HaLJiA4CCD1XsLUL7P628PHiaWR93cn8kXvmZKpxugWklhxzMUdAxO8NonF9N-uGQB4B7-4lX_3SD91JPzXuOLiPvsGDu3Menod9Pg6gkgWDHxOuh2ZRNrm5yFve3aLzysRcmg5diMQ92Yws4UuWjJvEXvD77q8MLTFxrgWAXOErjvtLv-E_gwG34fWMeduft6MSmf5F-oh8goqjjolGhi-h-dQQ5YecsoFCB3ncrJys.9EKnkXvyxCPIKZXs9Mx9U7m2t-LMTsbMy9fr4BxJK6M&dib_tag=se&keywords=gravitrax+junior+starterset&qid=1729516745&sprefix=gravitrax+junior%2Caps%2C98&sr=8-5"
$ws.Range("B25").Value = "https://m.media-amazon.com/images/I/81yqCdPH3SL._AC_SX679_.jpg"
$ws.Range("A25").Value = "Gravitrax Junior Starter Set Frozen"
$ws.Range("D25").Value = "36 EUR"

# Row 21 - Telefon Gabby
$ws.Range("C21").Value = "https://www.amazon.de/-/en/VTech-Gabbys-Dollhouse-Learning-Mobile/dp/B0CKRQ3N9W/ref=sr_1_2?crid=1TGUUBO2GUPXG&dib=eyJ2IjoiMSJ9.y1-F9U8wTKploOwLbra7x_dsi6XA_61WOmjL6YdtfKrMZd5PgbtYetS_xvrxVlN5P2N30MVC6N3-x3f18WvJLvsACyrQPzED_H1v5AXbB5G5syHB0uEk7crueAZVtpCq1DlsW_jJNYq258-TK7XttLHSAv_SDIynj2wHjEpJ8M9MAHFKYdOwMS1MVOfpFPOYGVf7604XryxZFTUkz68xV3oE47bLD-NFwuB1HEwlYbg.TibLhJEUOcQ47hE3d_W_l7dd4y7OSJytlm989H5OdP4&dib_tag=se&keywords=gabby+lernhandy&qid=1729517077&sprefix=gabby+lern%2Caps%2C108&sr=8-2#customerReviews"
$ws.Range("D21").Value = "20 EUR"
$ws.Range("B21").Value = "https://m.media-amazon.com/images/I/71L5GUCfx4L._AC_SX679_.jpg"
$ws.Range("A21").Value = "Telefon Gabby"

# Update selection / active cell to A22
[void]$ws.Range("A22").Select()
